# Convert the "Project List" sheet from its two-row merged-header layout
# (Group / Supervisor / Co-supervisor / Student spanning row 1, with the
# real column captions on row 2) into a flat, single-row-header table with
# camelCase column names, matching the shape used by the "readingexcel"
# sample code. Also add a new "AssessmentStatus" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old super-header row (old row 1: Group/Supervisor/Co-supervisor/Student).
#    This shifts every row up by one and collapses the merges that lived in
#    that row (A1:B1, C1:C2, D1:D2, E1:G1) automatically.
$ws.Rows.Item(1).Delete()

# 2) Rename the (now) row-1 header captions to the new camelCase keys.
$ws.Range("A1").Value = "groupNo"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "supervisor"
$ws.Range("D1").Value = "coSupervisor"
$ws.Range("E1").Value = "id"
$ws.Range("F1").Value = "lastName"
$ws.Range("G1").Value = "firstName"

# Bold + boxed styling to match the look of the rest of the header row.
$hdr = $ws.Range("A1:G1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.Borders.Item(7).LineStyle = 1
$hdr.Borders.Item(8).LineStyle = 1
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(10).LineStyle = 1
$hdr.Borders.Item(11).LineStyle = 1
$hdr.Borders.Item(12).LineStyle = 1

# The supervisor / co-supervisor header cells lose their bottom border so
# they read as a single banded header together with the data below.
$ws.Range("C1:D1").Borders.Item(9).LineStyle = 0

# 3) Add the new "AssessmentStatus" column (H) with a default value of 0
#    for every data row.
$ws.Range("H1").Value = "AssessmentStatus"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4108
$ws.Range("H1").Borders.Item(7).LineStyle = 1
$ws.Range("H1").Borders.Item(10).LineStyle = 1

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0

# 4) Match the saved selection/cursor position from the source edit.
$ws.Range("H7").Select() | Out-Null
